# The sheet contains one row per day of Ajo (garlic) price observations at
# "Terminal Hortofrutícola Agro Chillán". This edit inserts one additional
# daily observation before the existing row 26, shifting every subsequent
# row down by one (old row 26 -> new row 27, ..., old row 120 -> new row 121).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 26; everything from 26 downward shifts
# down by one row (this also grows the sheet dimension from R120 to R121).
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new observation. All
# non-numeric attributes (mercado, región, categoría, variedad, calidad,
# unidad, origen, clasificación) mirror the surrounding "Chino/Primera"
# Ajo entries; only the date, volume and price columns are new.
$ws.Range("A26").Value = 7
$ws.Range("B26").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C26").Value = "Ñuble"
$ws.Range("D26").Value = 44453
$ws.Range("E26").Value = 16
$ws.Range("F26").Value = 100112003
$ws.Range("G26").Value = "Ajo"
$ws.Range("H26").Value = "Chino"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 60
$ws.Range("K26").Value = 15000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 15500
$ws.Range("N26").Value = "$/caja 10 kilos"
$ws.Range("O26").Value = "China"
$ws.Range("P26").Value = 1550
$ws.Range("Q26").Value = 10
$ws.Range("R26").Value = "Hortaliza"
